# fix the work_mode initial value
#
# Sheet1 row 8 describes the "work_mode" variable (column B) but its
# "说明" (description) cell in column D still said "雷达开关" (radar
# switch) — a leftover/incorrect label. Correct it to "模块工作模式"
# (module work mode) so the description matches the variable it documents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D8").Value = "模块工作模式"

# Move the visible selection along with the edit.
$ws.Range("C11").Select()
